$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$ws.Range("A4").Value = "multiply"
$ws.Range("A5").Value = "divide"
$ws.Range("B4").Value = "Multiply"
$ws.Range("B5").Value = "Divide"

$ws.Range("B5").Select()
